# Fruta / hortaliza, semanal
# Insert two new weekly price rows into the "Frutilla" dataset right above the
# existing block that starts at row 357, pushing all the following rows down
# by two (so the former row 357 becomes row 359, etc.), and fill the two new
# rows with this week's values (copying the constant/template columns from
# the rows that end up right below them, then overwriting the columns that
# actually vary week to week: Fecha, Volumen, Precio minimo/maximo/promedio
# and Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 357.. down by two, creating two blank rows at 357 and 358.
$ws.Rows("357:358").Insert()

# New row 357 ("Especial" quality) - seed it from the row that is now right
# below it (old row 357, now shifted to 359) so all the constant columns
# (Mercado, Region, Producto, Calidad, Unidad, Origen, etc.) line up.
$ws.Range("A359:T359").Copy()
$ws.Range("A357:T357").PasteSpecial()

# New row 358 ("Primera" quality) - seed it the same way from row 360.
$ws.Range("A360:T360").Copy()
$ws.Range("A358:T358").PasteSpecial()

# Now overwrite the week-specific values for the two new rows.
$ws.Range("D357").Value = 44798
$ws.Range("M357").Value = 240
$ws.Range("N357").Value = 29000
$ws.Range("O357").Value = 30000
$ws.Range("P357").Value = 29500
$ws.Range("S357").Value = 4214

$ws.Range("D358").Value = 44798
$ws.Range("M358").Value = 160
$ws.Range("N358").Value = 26000
$ws.Range("O358").Value = 27000
$ws.Range("P358").Value = 26500
$ws.Range("S358").Value = 3786

$excel.CutCopyMode = $false
